$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Аркуш2")
$ws.Activate()

# --- "location" entity table (rows 18-24 in column I:L) -----------------
# The table used to start at row 19 with a spurious "Юридична назва"
# (legal name) field copy-pasted from the Company table; that field is
# removed, a new "email" field is appended at the end of the table, and
# the whole field block is renumbered down by one row.

# New first row of the table block: "email" field, right after "телефон"
# (row 17), replacing the previously blank spacer row 18.
$ws.Range("L18").Value = "email"
$ws.Rows.Item(18).AutoFit()

# Row 19 becomes the new (taller) blank spacer row that used to be row 18.
$ws.Range("I19:L19").Clear()
$ws.Rows.Item(19).RowHeight = 15.6

# Row 20: table header (owner_id / Company_id / location / TM), previously
# on row 19.
$ws.Range("D22").Copy()
$ws.Range("K20").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("I20").Value = "owner_id"
$ws.Range("J20").Value = "Company_id"
$ws.Range("K20").Value = "location"
$ws.Range("L20").Value = "TM"
$ws.Rows.Item(20).RowHeight = 15

# Row 21 ("Фактична адреса adress" + country/city/street/house/location)
# is unchanged.

# Rows 22-24: only the L-column field label shifts up one slot in the
# field list (ЄДРПОУ / телефон / фото drop out, телефон / фото / email
# move up); the B:E helper columns on these rows are untouched.
$ws.Range("L22").Value = "телефон"
$ws.Range("L23").Value = "фото"
$ws.Range("L24").Value = "email"

# --- view / selection ----------------------------------------------------
$ws.Range("N17").Select()
